# Add the function get_unique
# Normalizes header/category labels by replacing spaces with underscores
# and fixes the "Stuary" typo to "Etuary".
#
# Applying the edits cell-by-cell (rather than a blanket Find/Replace)
# mirrors how the values were actually retyped in the sheet, so the
# shared-string table is rebuilt in the same append order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Bin_name"
$ws.Range("B1").Value = "Sample_site"

$ws.Range("B4").Value = "Etuary"
$ws.Range("B5").Value = "Etuary"

$ws.Range("B2").Value = "Water_column"
$ws.Range("B7").Value = "Water_column"

$ws.Range("C2").Value = "Clade_1"
$ws.Range("C5").Value = "Clade_1"

$ws.Range("C3").Value = "Clade_2"
$ws.Range("C6").Value = "Clade_2"

$ws.Range("C4").Value = "Clade_3"
$ws.Range("C7").Value = "Clade_3"

# Widen the first two columns now that they hold longer, underscored labels
$ws.Columns.Item(1).ColumnWidth = 31.330729166666668
$ws.Columns.Item(2).ColumnWidth = 34.330729166666664

# Update the active cell selection
[void]$ws.Range("G12").Select()
